$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
}

Set-TextValue $ws.Range("D2") '37.288.31'
Set-TextValue $ws.Range("E2") '  -1.34%  '
Set-TextValue $ws.Range("D3") '2.047.08'
Set-TextValue $ws.Range("E3") '  -1.45%  '
Set-TextValue $ws.Range("E4") '  -0.07%  '
Set-TextValue $ws.Range("D5") '228.95'
Set-TextValue $ws.Range("E5") '  -1.67%  '
Set-TextValue $ws.Range("D6") '0.614'
Set-TextValue $ws.Range("E6") '  -1.73%  '
Set-TextValue $ws.Range("E7") '  +0.04%  '
Set-TextValue $ws.Range("D8") '56.47'
Set-TextValue $ws.Range("E8") '  -3.43%  '
Set-TextValue $ws.Range("E9") '  -2.07%  '
Set-TextValue $ws.Range("D10") '0.0786'
Set-TextValue $ws.Range("E10") '  +0.11%  '
Set-TextValue $ws.Range("E11") '  -2.02%  '
Set-TextValue $ws.Range("D12") '14.69'
Set-TextValue $ws.Range("E12") '  -1.36%  '
Set-TextValue $ws.Range("D13") '2.328.73'
Set-TextValue $ws.Range("E13") '  -2.31%  '
Set-TextValue $ws.Range("D14") '20.60'
Set-TextValue $ws.Range("E14") '  -2.17%  '
Set-TextValue $ws.Range("E15") '  -3.27%  '
Set-TextValue $ws.Range("E16") '  -1.04%  '
Set-TextValue $ws.Range("D17") '2.048.26'
Set-TextValue $ws.Range("E17") '  -0.50%  '
Set-TextValue $ws.Range("D18") '37.197.68'
Set-TextValue $ws.Range("E18") '  -1.42%  '
Set-TextValue $ws.Range("D19") '6.07'
Set-TextValue $ws.Range("E19") '  -0.88%  '
Set-TextValue $ws.Range("D20") '69.36'
Set-TextValue $ws.Range("E20") '  -3.09%  '
Set-TextValue $ws.Range("E21") '  -2.00%  '
Set-TextValue $ws.Range("D22") '225.56'
Set-TextValue $ws.Range("E22") '  -1.54%  '
Set-TextValue $ws.Range("D23") '1.00'
Set-TextValue $ws.Range("E23") '  +0.11%  '
Set-TextValue $ws.Range("D24") '2.40'
Set-TextValue $ws.Range("E24") '  +0.28%  '
Set-TextValue $ws.Range("D25") '2.28'
Set-TextValue $ws.Range("E25") '  -4.86%  '
Set-TextValue $ws.Range("D26") '9.68'
Set-TextValue $ws.Range("E26") '  +0.44%  '
Set-TextValue $ws.Range("D27") '167.12'
Set-TextValue $ws.Range("E27") '  -2.70%  '
Set-TextValue $ws.Range("E28") '  -7.19%  '
Set-TextValue $ws.Range("E29") '  -2.37%  '
Set-TextValue $ws.Range("E30") '  -4.19%  '
Set-TextValue $ws.Range("E31") '  -1.88%  '
Set-TextValue $ws.Range("D32") '4.53'
Set-TextValue $ws.Range("E32") '  -4.26%  '
Set-TextValue $ws.Range("E33") '  -2.86%  '
Set-TextValue $ws.Range("D34") '4.58'
Set-TextValue $ws.Range("E34") '  -1.82%  '
Set-TextValue $ws.Range("E35") '  -1.18%  '
Set-TextValue $ws.Range("E36") '  +1.31%  '
Set-TextValue $ws.Range("E37") '  +0.09%  '
Set-TextValue $ws.Range("E38") '  -4.94%  '
Set-TextValue $ws.Range("D39") '5.23'
Set-TextValue $ws.Range("E39") '  -3.08%  '
Set-TextValue $ws.Range("E40") '  -4.59%  '
Set-TextValue $ws.Range("E41") '  -1.03%  '
Set-TextValue $ws.Range("D42") '1.476.81'
Set-TextValue $ws.Range("E42") '  +1.72%  '
Set-TextValue $ws.Range("D43") '16.93'
Set-TextValue $ws.Range("E43") '  -0.48%  '
Set-TextValue $ws.Range("D44") '96.29'
Set-TextValue $ws.Range("E44") '  -5.95%  '
Set-TextValue $ws.Range("E45") '  -3.76%  '
Set-TextValue $ws.Range("E46") '  +0.73%  '
Set-TextValue $ws.Range("E47") '  -4.29%  '
Set-TextValue $ws.Range("D48") '3.92'
Set-TextValue $ws.Range("E48") '  -4.59%  '
Set-TextValue $ws.Range("D49") '7.08'
Set-TextValue $ws.Range("E49") '  -3.77%  '
Set-TextValue $ws.Range("E50") '  -2.30%  '
Set-TextValue $ws.Range("D51") '2.232.37'
Set-TextValue $ws.Range("E51") '  -1.64%  '
